# Tabla de bisección: añadir tipo de error en sección 1, finalización sistema y radio sor
# - Actualiza la columna E (error) de las filas 3 a 21 para reflejar el error
#   relativo (|xn - xn-1| / |xn|) en lugar del error absoluto.
# - Elimina la última fila (22), que ya no corresponde al nuevo criterio de
#   finalización del sistema.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newErrors = @{
  3  = "1.0"
  4  = "0.333333333333333"
  5  = "0.2"
  6  = "0.0909090909090909"
  7  = "0.0476190476190476"
  8  = "0.024390243902439"
  9  = "0.0123456790123457"
  10 = "0.0061349693251533"
  11 = "0.003076923076923"
  12 = "0.0015408320493066"
  13 = "0.0007710100231303"
  14 = "0.0003856536829926"
  15 = "0.0001928640308582"
  16 = "9.64413154595429e-05"
  17 = "4.82229830737329e-05"
  18 = "2.41109101868596e-05"
  19 = "1.20556004291794e-05"
  20 = "6.02783654918413e-06"
  21 = "3.01390919091608e-06"
}

foreach ($r in $newErrors.Keys) {
  $cell = $ws.Cells.Item($r, 4)
  # Forzar el valor como texto (para conservar el formato de cadena original)
  # y luego restablecer el estilo por defecto para no dejar un formato de
  # número/texto aplicado a la celda.
  $cell.NumberFormat = "@"
  $cell.Value = $newErrors[$r]
  $cell.Style = "Normal"
}

# Elimina la fila 22 (última iteración), que ya no aparece en la tabla.
$ws.Range("A22:D22").EntireRow.Delete()
